# Add a new "Sheet2" data-driven test-data sheet after the existing "Sheet1",
# matching the "data writing to excel code is added" commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so it keeps sheetId/rId ordering
# (Sheet1 stays sheetId=1/rId1, Sheet2 becomes sheetId=2/rId2) and becomes the
# active tab, same as Excel does when you right-click > Insert a sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Helper: write a value as literal text (not auto-coerced to a number) with
# no left-over cell styling, for numeric-looking ids like "11538380".
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Data rows entered first ...
Set-TextValue $ws2.Range("A2") "TN2485269"
Set-TextValue $ws2.Range("B2") "11538380"

Set-TextValue $ws2.Range("A3") "TB2485272"
Set-TextValue $ws2.Range("B3") "11538402"

# ... then the header row ...
$ws2.Range("A1").Value = "Policy Number"
$ws2.Range("B1").Value = "Quote Number"
$ws2.Range("C1").Value = "Type Policy"
$ws2.Range("D1").Value = "Type"
$ws2.Range("E1").Value = "Testcaseid"

# ... then a final data row.
Set-TextValue $ws2.Range("A4") "TB2485273"
Set-TextValue $ws2.Range("B4") "11538440"
$ws2.Range("C4").Value = "Bond - No Credit"
$ws2.Range("D4").Value = "TC001"

# Leave the cursor where the author left it before saving.
$ws2.Activate()
$ws2.Range("I8").Select()
